# Auto-generated Excel COM-interop script
# Applies targeted cell-value corrections across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets
# (WVR untouched), matching the scheduled-runner price-refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 224.4
$ws.Range("I12").Value2 = 199.66667
$ws.Range("K12").Value2 = 199.66667
$ws.Range("M12").Value2 = -29.66667000000001

$ws.Range("H48").Value2 = 500
$ws.Range("J48").Value2 = 500
$ws.Range("L48").Value2 = 1500
$ws.Range("N48").Value2 = -2084

$ws.Range("H56").Value2 = 500
$ws.Range("J56").Value2 = 500
$ws.Range("L56").Value2 = 1500
$ws.Range("N56").Value2 = -2568

$ws.Range("H86").Value2 = 2999.25
$ws.Range("I86").Value2 = 2999.25
$ws.Range("K86").Value2 = 2999.25
$ws.Range("M86").Value2 = -1876.25

$ws.Range("H89").Value2 = 2999.25
$ws.Range("I89").Value2 = 2999.25
$ws.Range("K89").Value2 = 14996.25
$ws.Range("M89").Value2 = -9380.25

$ws.Range("H111").Value2 = 0
$ws.Range("I111").Value2 = 0
$ws.Range("K111").Value2 = 0
$ws.Range("M111").Value2 = ""

$ws.Range("H125").Value2 = 849.6667
$ws.Range("I125").Value2 = 774.5
$ws.Range("J125").Value2 = 1000
$ws.Range("K125").Value2 = 6970.5
$ws.Range("L125").Value2 = 9000
$ws.Range("M125").Value2 = -4510.5
$ws.Range("N125").Value2 = -13920

$ws.Range("H137").Value2 = 3081.818
$ws.Range("I137").Value2 = 2783.3333
$ws.Range("K137").Value2 = 8349.999899999999
$ws.Range("M137").Value2 = -5799.999899999999

$ws.Range("H141").Value2 = 875
$ws.Range("I141").Value2 = 875
$ws.Range("K141").Value2 = 2625
$ws.Range("M141").Value2 = 2555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 7868.1816
$ws.Range("I32").Value2 = 7905
$ws.Range("J32").Value2 = 7500
$ws.Range("K32").Value2 = 7905
$ws.Range("L32").Value2 = 7500
$ws.Range("M32").Value2 = -7618
$ws.Range("N32").Value2 = -8074

$ws.Range("H102").Value2 = 2480
$ws.Range("I102").Value2 = 2480
$ws.Range("K102").Value2 = 2480
$ws.Range("M102").Value2 = -858

$ws.Range("H122").Value2 = 2101.7144
$ws.Range("J122").Value2 = 3500
$ws.Range("L122").Value2 = 10500
$ws.Range("N122").Value2 = -15400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 2683
$ws.Range("I20").Value2 = 2567.6667
$ws.Range("J20").Value2 = 2798.3333
$ws.Range("K20").Value2 = 2567.6667
$ws.Range("L20").Value2 = 2798.3333
$ws.Range("M20").Value2 = -2320.6667
$ws.Range("N20").Value2 = -3292.3333

$ws.Range("H81").Value2 = 39890
$ws.Range("J81").Value2 = 39890
$ws.Range("L81").Value2 = 39890
$ws.Range("N81").Value2 = -42012

$ws.Range("H84").Value2 = 39890
$ws.Range("J84").Value2 = 39890
$ws.Range("L84").Value2 = 119670
$ws.Range("N84").Value2 = -130278

$ws.Range("H86").Value2 = 0
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 0
$ws.Range("L86").Value2 = 0
$ws.Range("M86").Value2 = ""
$ws.Range("N86").Value2 = ""

$ws.Range("H89").Value2 = 0
$ws.Range("I89").Value2 = 0
$ws.Range("J89").Value2 = 0
$ws.Range("K89").Value2 = 0
$ws.Range("L89").Value2 = 0
$ws.Range("M89").Value2 = ""
$ws.Range("N89").Value2 = ""

$ws.Range("H105").Value2 = 150000
$ws.Range("I105").Value2 = 150000
$ws.Range("J105").Value2 = 0
$ws.Range("K105").Value2 = 150000
$ws.Range("L105").Value2 = 0
$ws.Range("M105").Value2 = -148253
$ws.Range("N105").Value2 = ""

$ws.Range("H107").Value2 = 1398.5
$ws.Range("I107").Value2 = 1464.6666
$ws.Range("K107").Value2 = 1464.6666
$ws.Range("M107").Value2 = 455.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 447.5
$ws.Range("I22").Value2 = 447.5
$ws.Range("K22").Value2 = 447.5
$ws.Range("M22").Value2 = -97.5

$ws.Range("H31").Value2 = 2000
$ws.Range("I31").Value2 = 2000
$ws.Range("K31").Value2 = 2000
$ws.Range("M31").Value2 = -1705

$ws.Range("H34").Value2 = 2000
$ws.Range("I34").Value2 = 2000
$ws.Range("K34").Value2 = 2000
$ws.Range("M34").Value2 = -1798

$ws.Range("H86").Value2 = 0
$ws.Range("I86").Value2 = 0
$ws.Range("K86").Value2 = 0
$ws.Range("M86").Value2 = ""

$ws.Range("H89").Value2 = 0
$ws.Range("I89").Value2 = 0
$ws.Range("K89").Value2 = 0
$ws.Range("M89").Value2 = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value2 = 2995
$ws.Range("J31").Value2 = 2995
$ws.Range("L31").Value2 = 8985
$ws.Range("N31").Value2 = -9561

$ws.Range("H38").Value2 = 1577.7273
$ws.Range("I38").Value2 = 1344
$ws.Range("J38").Value2 = 1858.2
$ws.Range("K38").Value2 = 4032
$ws.Range("L38").Value2 = 5574.6
$ws.Range("M38").Value2 = -3685
$ws.Range("N38").Value2 = -6268.6

$ws.Range("H40").Value2 = 90
$ws.Range("I40").Value2 = 75
$ws.Range("K40").Value2 = 300
$ws.Range("M40").Value2 = -231

$ws.Range("H46").Value2 = 727
$ws.Range("I46").Value2 = 500
$ws.Range("K46").Value2 = 1500
$ws.Range("M46").Value2 = -1409

$ws.Range("H56").Value2 = 0
$ws.Range("I56").Value2 = 0
$ws.Range("K56").Value2 = 0
$ws.Range("M56").Value2 = ""

$ws.Range("H80").Value2 = 5000
$ws.Range("I80").Value2 = 5000
$ws.Range("K80").Value2 = 15000
$ws.Range("M80").Value2 = -14064

$ws.Range("H83").Value2 = 5000
$ws.Range("I83").Value2 = 5000
$ws.Range("K83").Value2 = 45000
$ws.Range("M83").Value2 = -40320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2000
$ws.Range("I80").Value2 = 1000
$ws.Range("J80").Value2 = 3000
$ws.Range("K80").Value2 = 1000
$ws.Range("L80").Value2 = 3000
$ws.Range("M80").Value2 = -2
$ws.Range("N80").Value2 = -4996

$ws.Range("H83").Value2 = 2000
$ws.Range("I83").Value2 = 1000
$ws.Range("J83").Value2 = 3000
$ws.Range("K83").Value2 = 5000
$ws.Range("L83").Value2 = 15000
$ws.Range("M83").Value2 = -8
$ws.Range("N83").Value2 = -24984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value2 = 3363.3333
$ws.Range("I2").Value2 = 10000
$ws.Range("J2").Value2 = 45
$ws.Range("K2").Value2 = 10000
$ws.Range("L2").Value2 = 45
$ws.Range("M2").Value2 = -9888
$ws.Range("N2").Value2 = -269

$ws.Range("H22").Value2 = 0
$ws.Range("I22").Value2 = 0
$ws.Range("K22").Value2 = 0
$ws.Range("M22").Value2 = ""

$ws.Range("H27").Value2 = 0
$ws.Range("I27").Value2 = 0
$ws.Range("K27").Value2 = 0
$ws.Range("M27").Value2 = ""

$ws.Range("H61").Value2 = 2514.4443
$ws.Range("I61").Value2 = 2616
$ws.Range("J61").Value2 = 2387.5
$ws.Range("K61").Value2 = 2616
$ws.Range("L61").Value2 = 2387.5
$ws.Range("M61").Value2 = -2414
$ws.Range("N61").Value2 = -2791.5

$ws.Range("H62").Value2 = 7437.25
$ws.Range("J62").Value2 = 9374.5
$ws.Range("L62").Value2 = 9374.5
$ws.Range("N62").Value2 = -10622.5

$ws.Range("H65").Value2 = 7437.25
$ws.Range("J65").Value2 = 9374.5
$ws.Range("L65").Value2 = 28123.5
$ws.Range("N65").Value2 = -34363.5

$ws.Range("H93").Value2 = 600.2
$ws.Range("I93").Value2 = 500.375
$ws.Range("J93").Value2 = 999.5
$ws.Range("K93").Value2 = 500.375
$ws.Range("L93").Value2 = 999.5
$ws.Range("M93").Value2 = 747.625
$ws.Range("N93").Value2 = -3495.5

$ws.Range("H113").Value2 = 2514.4443
$ws.Range("I113").Value2 = 2616
$ws.Range("J113").Value2 = 2387.5
$ws.Range("K113").Value2 = 2616
$ws.Range("L113").Value2 = 2387.5
$ws.Range("M113").Value2 = -446
$ws.Range("N113").Value2 = -6727.5

$ws.Range("H136").Value2 = 28336.666
$ws.Range("I136").Value2 = 10004
$ws.Range("K136").Value2 = 30012
$ws.Range("M136").Value2 = -27462
